# edit.ps1 -- apply the Analisi_requisiti.docx changes described in the diff:
#  1) RD3 paragraph: split "RD3:" / "Un" runs and add proofErr gramStart/gramEnd
#  2) RD5 (menu fissi) paragraph: split "primo,secondo" out and add proofErr gramStart/gramEnd
#  3) RD6 paragraph: split "RD6:" / "Il" runs and add proofErr gramStart/gramEnd
#  4) RF7 paragraph: reword "le quantità degli ingredienti richiesti dal magazzino"
#     -> ", precise quantità di ingredienti dal magazzino"

$d = $word.ActiveDocument

function Get-ParagraphRangeByAnchor($anchorText) {
    $search = $d.Content
    $found = $search.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $anchorText"
    }
    $para = $search.Paragraphs(1)
    $pr = $para.Range
    return $d.Range($pr.Start, $pr.End)
}

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1) RD3 paragraph ---------------------------------------------------
$rd3 = Get-ParagraphRangeByAnchor("RD3:Un ordine contiene il numero di posti effettivamente occupati nel tavolo e tutti i piatti e menu ordinati dai clienti.")
$rd3xml = $pkgHeader + '<w:p w14:paraId="79BD6740" w14:textId="5149297C" w:rsidR="0060578C" w:rsidRDefault="0060578C" w:rsidP="0060578C"><w:r><w:t>RD</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>3:</w:t></w:r><w:r><w:t>Un</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ordine contiene il numero di posti effettivamente occupati nel tavolo e tutti i piatti e menu ordinati dai clienti.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' + $pkgFooter
$rd3.InsertXML($rd3xml)

# --- 2) RD5 (menu fissi) paragraph --------------------------------------
$rd5 = Get-ParagraphRangeByAnchor("RD5: I menu fissi sono formati da primo,secondo,frutta, dolce ed eventualmente un contorno.")
$rd5xml = $pkgHeader + '<w:p w14:paraId="0B1A2EDC" w14:textId="595770CA" w:rsidR="0091142C" w:rsidRDefault="0091142C" w:rsidP="0060578C"><w:r><w:t xml:space="preserve">RD5: I menu fissi sono formati da </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>primo,secondo</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>,frutta, dolce ed eventualmente un contorno.</w:t></w:r></w:p>' + $pkgFooter
$rd5.InsertXML($rd5xml)

# --- 3) RD6 paragraph ----------------------------------------------------
$rd6 = Get-ParagraphRangeByAnchor("RD6:Il menu memorizza tutte le pietanze ordinabili.")
$rd6xml = $pkgHeader + '<w:p w14:paraId="3A2A2EBF" w14:textId="4A903514" w:rsidR="0060578C" w:rsidRDefault="0060578C" w:rsidP="0060578C"><w:r><w:t>RD</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>6:</w:t></w:r><w:r><w:t>Il</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> menu memorizza tutte le pietanze ordinabili.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' + $pkgFooter
$rd6.InsertXML($rd6xml)

# --- 4) RF7 paragraph (prenotare / magazzino rewording) ------------------
$rf7 = Get-ParagraphRangeByAnchor("RF7: Il sistema deve fornire una funzionalità per prenotare le quantità degli ingredienti richiesti dal magazzino.")
$rf7xml = $pkgHeader + '<w:p w14:paraId="44588E75" w14:textId="7DD9A845" w:rsidR="00B5256B" w:rsidRDefault="00B5256B" w:rsidP="00B5256B"><w:r><w:t>RF</w:t></w:r><w:r w:rsidR="000A1613"><w:t>7</w:t></w:r><w:r w:rsidR="00F735E3"><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> Il sistema deve</w:t></w:r><w:r w:rsidR="000A1613"><w:t xml:space="preserve"> fornire una funzionalità per</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="000A1613"><w:t>prenotare</w:t></w:r><w:r><w:t>, precise quantità di ingredienti</w:t></w:r><w:r><w:t xml:space="preserve"> dal magazzino</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>' + $pkgFooter
$rf7.InsertXML($rf7xml)

Write-Output "Done."
